# "Update Work Week and Social Spending"
#
# The "Data" sheet holds one GDP-per-capita series for Azerbaijan (Country
# Code 31), one row per year, in column E ("Data"), stored as text (the
# source feed writes numeric-looking values as strings). This update:
#   1. Refreshes the GDP-per-capita figure for 1973 (row 2).
#   2. Fills in values for 1980-1989 (rows 9-18), which were previously
#      blank, and refreshes 1990-2010 (rows 19-39) with revised figures.
#   3. Appends six new years, 2011-2016 (rows 40-45).
#
# Column E values are written as text (not numbers) to match the source
# column's storage convention, so NumberFormat is forced to "@" before each
# assignment -- otherwise the host would auto-coerce a numeric-looking
# string like "7165" into a number cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- 1. Row 2 (year 1973): refreshed GDP-per-capita figure ---------------
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "7068"

# --- 2. Rows 9-39 (years 1980-2010): new/refreshed figures ---------------
$existingRows = @(
    @(9, "7165"), @(10, "7490"), @(11, "7831"), @(12, "8010"),
    @(13, "8360"), @(14, "8558"), @(15, "8663"), @(16, "8995"),
    @(17, "9328"), @(18, "8421"), @(19, "7394"),
    @(20, "7319.13646523413"), @(21, "5638.99790682118"),
    @(22, "4315.35671425961"), @(23, "3453.36382182429"),
    @(24, "2999.31904512599"), @(25, "3074.21973457034"),
    @(26, "3347.61264566942"), @(27, "3550.36406585758"),
    @(28, "3959.78147692477"), @(29, "4214.88886074828"),
    @(30, "4498.79437699741"), @(31, "4880.39989008972"),
    @(32, "5408.52298962295"), @(33, "5971.47460310057"),
    @(34, "7553.62639595617"), @(35, "10160.3460024994"),
    @(36, "12699.4977767115"), @(37, "14072.6265772429"),
    @(38, "15384.2930387916"), @(39, "16153.8365683861")
)

foreach ($item in $existingRows) {
    $row = $item[0]
    $val = $item[1]
    $cell = $ws.Range("E" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# --- 3. Rows 40-45 (years 2011-2016): brand-new rows ----------------------
$newRows = @(
    @(40, 2011, "16176"),
    @(41, 2012, "16359"),
    @(42, 2013, "17133"),
    @(43, 2014, "17439"),
    @(44, 2015, "17460"),
    @(45, 2016, "16645")
)

foreach ($item in $newRows) {
    $row = $item[0]
    $year = $item[1]
    $val = $item[2]

    $ws.Range("A" + $row).Value = 31
    $ws.Range("B" + $row).Value = "Azerbaijan"
    $ws.Range("C" + $row).Value = "GDP per Capita"
    $ws.Range("D" + $row).Value = $year

    $eCell = $ws.Range("E" + $row)
    $eCell.NumberFormat = "@"
    $eCell.Value = $val
}
